# Apply edit: rename the "Job Search.xlsx" module file reference (removing the
# space) and mark the cell with the built-in "Hyperlink" style, then move the
# active selection to D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the displayed text of D4 (keeps its existing hyperlink target).
$ws.Range("D4").Value = "\\Resources\\ModuleFiles\\JobSearch.xlsx"

# Apply the built-in Hyperlink style (adds Hyperlink font/cellStyle entries).
$ws.Range("D4").Style = "Hyperlink"

# Update the sheet's active selection.
$ws.Range("D4").Select()
